# Auto-generated Excel COM-interop script
# Applies scheduled market-price updates to the Leve profit sheets
# (columns H-N: currentAveragePrice.. through LeveProfitHQ) across all 8 job sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 25399.5
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

$ws.Range("H28").Value = 220
$ws.Range("J28").Value = 800
$ws.Range("L28").Value = 800
$ws.Range("N28").Value = -1770

$ws.Range("H32").Value = 997.5
$ws.Range("I32").Value = 997.5
$ws.Range("K32").Value = 997.5
$ws.Range("M32").Value = -671.5

$ws.Range("H33").Value = 200.94737
$ws.Range("I33").Value = 210.17647
$ws.Range("K33").Value = 210.17647
$ws.Range("M33").Value = 18.82353000000001

$ws.Range("H40").Value = 2461.3333
$ws.Range("I40").Value = 2461.3333
$ws.Range("K40").Value = 2461.3333
$ws.Range("M40").Value = -2286.3333

$ws.Range("H112").Value = 1466.6666
$ws.Range("J112").Value = 1466.6666
$ws.Range("L112").Value = 4399.9998
$ws.Range("N112").Value = -6615.9998

$ws.Range("H129").Value = 1941.3334
$ws.Range("J129").Value = 2517
$ws.Range("L129").Value = 7551
$ws.Range("N129").Value = -17551

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3069.7778
$ws.Range("I74").Value = 2266.6
$ws.Range("J74").Value = 4073.75
$ws.Range("K74").Value = 2266.6
$ws.Range("L74").Value = 4073.75
$ws.Range("M74").Value = -1392.6
$ws.Range("N74").Value = -5821.75

$ws.Range("H77").Value = 3069.7778
$ws.Range("I77").Value = 2266.6
$ws.Range("J77").Value = 4073.75
$ws.Range("K77").Value = 11333
$ws.Range("L77").Value = 20368.75
$ws.Range("M77").Value = -6965
$ws.Range("N77").Value = -29104.75

$ws.Range("H103").Value = 10000
$ws.Range("J103").Value = 10000
$ws.Range("L103").Value = 10000
$ws.Range("N103").Value = -12344

$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

$ws.Range("H132").Value = 5635.625
$ws.Range("I132").Value = 4295
$ws.Range("K132").Value = 12885
$ws.Range("M132").Value = -10355

$ws.Range("H138").Value = 55000
$ws.Range("I138").Value = 55000
$ws.Range("K138").Value = 55000
$ws.Range("M138").Value = -49860

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 70000
$ws.Range("J50").Value = 70000
$ws.Range("L50").Value = 70000
$ws.Range("N50").Value = -71148

$ws.Range("H80").Value = 925.2
$ws.Range("J80").Value = 969.5
$ws.Range("L80").Value = 969.5
$ws.Range("N80").Value = -2965.5

$ws.Range("H83").Value = 925.2
$ws.Range("J83").Value = 969.5
$ws.Range("L83").Value = 4847.5
$ws.Range("N83").Value = -14831.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 2840
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

$ws.Range("H86").Value = 4735.6665
$ws.Range("I86").Value = 3103.5
$ws.Range("J86").Value = 8000
$ws.Range("K86").Value = 3103.5
$ws.Range("L86").Value = 8000
$ws.Range("M86").Value = -1980.5
$ws.Range("N86").Value = -10246

$ws.Range("H89").Value = 4735.6665
$ws.Range("I89").Value = 3103.5
$ws.Range("J89").Value = 8000
$ws.Range("K89").Value = 15517.5
$ws.Range("L89").Value = 40000
$ws.Range("M89").Value = -9901.5
$ws.Range("N89").Value = -51232

$ws.Range("H93").Value = 10000
$ws.Range("I93").Value = 10000
$ws.Range("K93").Value = 10000
$ws.Range("M93").Value = -8128

$ws.Range("H94").Value = 5833.3335
$ws.Range("I94").Value = 8250
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 8250
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = -7799
$ws.Range("N94").Value = -1902

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("N92").ClearContents()

$ws.Range("H97").Value = 480.66666
$ws.Range("I97").Value = 346
$ws.Range("J97").Value = 750
$ws.Range("K97").Value = 1038
$ws.Range("L97").Value = 2250
$ws.Range("M97").Value = -542
$ws.Range("N97").Value = -3242

$ws.Range("H131").Value = 2612.1333
$ws.Range("J131").Value = 2760.8572
$ws.Range("L131").Value = 8282.571599999999
$ws.Range("N131").Value = -18362.5716

$ws.Range("H132").Value = 500
$ws.Range("I132").Value = 500
$ws.Range("K132").Value = 4500
$ws.Range("M132").Value = -1970

$ws.Range("H133").Value = 3000
$ws.Range("J133").Value = 3000
$ws.Range("L133").Value = 9000
$ws.Range("N133").Value = -19120

$ws.Range("H134").Value = 4095
$ws.Range("J134").Value = 6000
$ws.Range("L134").Value = 18000
$ws.Range("N134").Value = -28140

$ws.Range("H139").Value = 2415.375
$ws.Range("I139").Value = 1774.3334
$ws.Range("J139").Value = 2800
$ws.Range("K139").Value = 5323.0002
$ws.Range("L139").Value = 8400
$ws.Range("M139").Value = -183.0002000000004
$ws.Range("N139").Value = -18680

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 39999
$ws.Range("J134").Value = 39999
$ws.Range("L134").Value = 119997
$ws.Range("N134").Value = -125067

$ws.Range("H141").Value = 66999.5
$ws.Range("J141").Value = 66999.5
$ws.Range("L141").Value = 66999.5
$ws.Range("N141").Value = -77359.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4056.6155
$ws.Range("I7").Value = 4089.6365
$ws.Range("K7").Value = 4089.6365
$ws.Range("M7").Value = -3977.6365

$ws.Range("H46").Value = 1131.4
$ws.Range("J46").Value = 1131.4
$ws.Range("L46").Value = 1131.4
$ws.Range("N46").Value = -1507.4

$ws.Range("H55").Value = 674.5714
$ws.Range("I55").Value = 620
$ws.Range("K55").Value = 620
$ws.Range("M55").Value = -447

$ws.Range("H101").Value = 14998
$ws.Range("J101").Value = 14998
$ws.Range("L101").Value = 14998
$ws.Range("N101").Value = -21488

$ws.Range("H103").Value = 20501
$ws.Range("J103").Value = 20501
$ws.Range("L103").Value = 20501
$ws.Range("N103").Value = -22845

$ws.Range("H126").Value = 4056.6155
$ws.Range("I126").Value = 4089.6365
$ws.Range("K126").Value = 12268.9095
$ws.Range("M126").Value = -9798.9095

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
